$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.307.97"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.922.21"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8126"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.52"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3261"
$ws.Range("E8").Value = "  +2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.25"
$ws.Range("E9").Value = "  +3.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07247"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7943"
$ws.Range("E11").Value = "  +6.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08120"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.917.97"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.445"
$ws.Range("E14").Value = "  +4.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.51"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.302.40"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.29"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.093"
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "250.20"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.182.20"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.250"
$ws.Range("E22").Value = "  +20.50%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1651"
$ws.Range("E25").Value = "  +18.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.532"
$ws.Range("E26").Value = "  +3.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.13"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.163"
$ws.Range("E29").Value = "  +6.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.373"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.554"
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.356"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05768"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.146"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.307"
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7494"
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9993"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.728"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01961"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.820"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4515"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.94"
$ws.Range("E42").Value = "  +3.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.996"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.934"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.042.80"
$ws.Range("E46").Value = "  +5.59%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.55"
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.116"
$ws.Range("E49").Value = "  +10.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.667"
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.923"
$ws.Range("E51").Value = "  +1.43%  "
